# Apply the 30-Mar-2020 23:50 COVID-19 data refresh to the "Pais" sheet.
# (Updated totals for several countries; a few rows swap position because
#  the sheet is kept sorted by total cases (column B) descending.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in cell A1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 30 de Marzo de 2020 a las 23:50"

# Row 4: Estados Unidos - refreshed totals
$ws.Cells.Item(4,1).Value = "Estados Unidos"
$ws.Cells.Item(4,2).Value = 160385
$ws.Cells.Item(4,3).Value = 16894
$ws.Cells.Item(4,4).Value = 5220
$ws.Cells.Item(4,5).Value = 152212
$ws.Cells.Item(4,6).Value = 3402
$ws.Cells.Item(4,7).Value = 370
$ws.Cells.Item(4,8).Value = 2953

# Row 6: España - refreshed totals
$ws.Cells.Item(6,1).Value = "España"
$ws.Cells.Item(6,2).Value = 87956
$ws.Cells.Item(6,3).Value = 7846
$ws.Cells.Item(6,4).Value = 16780
$ws.Cells.Item(6,5).Value = 63460
$ws.Cells.Item(6,6).Value = 5231
$ws.Cells.Item(6,7).Value = 913
$ws.Cells.Item(6,8).Value = 7716

# Row 8: Alemania - refreshed totals
$ws.Cells.Item(8,1).Value = "Alemania"
$ws.Cells.Item(8,2).Value = 66711
$ws.Cells.Item(8,3).Value = 4276
$ws.Cells.Item(8,4).Value = 13500
$ws.Cells.Item(8,5).Value = 52566
$ws.Cells.Item(8,6).Value = 1979
$ws.Cells.Item(8,7).Value = 104
$ws.Cells.Item(8,8).Value = 645

# Row 12: Suiza - refreshed totals
$ws.Cells.Item(12,1).Value = "Suiza"
$ws.Cells.Item(12,2).Value = 15922
$ws.Cells.Item(12,3).Value = 1093
$ws.Cells.Item(12,4).Value = 1823
$ws.Cells.Item(12,5).Value = 13740
$ws.Cells.Item(12,6).Value = 301
$ws.Cells.Item(12,7).Value = 59
$ws.Cells.Item(12,8).Value = 359

# Row 20: Israel - refreshed totals
$ws.Cells.Item(20,1).Value = "Israel"
$ws.Cells.Item(20,2).Value = 4695
$ws.Cells.Item(20,3).Value = 448
$ws.Cells.Item(20,4).Value = 161
$ws.Cells.Item(20,5).Value = 4518
$ws.Cells.Item(20,6).Value = 79
$ws.Cells.Item(20,7).Value = 1
$ws.Cells.Item(20,8).Value = 16

# Row 23: Australia - refreshed totals
$ws.Cells.Item(23,1).Value = "Australia"
$ws.Cells.Item(23,2).Value = 4359
$ws.Cells.Item(23,3).Value = 196
$ws.Cells.Item(23,4).Value = 244
$ws.Cells.Item(23,5).Value = 4097
$ws.Cells.Item(23,6).Value = 28
$ws.Cells.Item(23,7).Value = 1
$ws.Cells.Item(23,8).Value = 18

# Row 36: Pakistan - refreshed totals
$ws.Cells.Item(36,1).Value = "Pakistan"
$ws.Cells.Item(36,2).Value = 1717
$ws.Cells.Item(36,3).Value = 120
$ws.Cells.Item(36,4).Value = 76
$ws.Cells.Item(36,5).Value = 1620
$ws.Cells.Item(36,6).Value = 11
$ws.Cells.Item(36,7).Value = 7
$ws.Cells.Item(36,8).Value = 21

# Row 65: now "Marruecos" (was "Ucrania") - sheet re-sorted by total cases
$ws.Cells.Item(65,1).Value = "Marruecos"
$ws.Cells.Item(65,2).Value = 556
$ws.Cells.Item(65,3).Value = 77
$ws.Cells.Item(65,4).Value = 15
$ws.Cells.Item(65,5).Value = 508
$ws.Cells.Item(65,6).Value = 1
$ws.Cells.Item(65,7).Value = 7
$ws.Cells.Item(65,8).Value = 33

# Row 66: now "Ucrania" (was "Marruecos") - sheet re-sorted by total cases
$ws.Cells.Item(66,1).Value = "Ucrania"
$ws.Cells.Item(66,2).Value = 548
$ws.Cells.Item(66,3).Value = 73
$ws.Cells.Item(66,4).Value = 8
$ws.Cells.Item(66,5).Value = 527
$ws.Cells.Item(66,6).Value = 0
$ws.Cells.Item(66,7).Value = 3
$ws.Cells.Item(66,8).Value = 13

# Row 103: Honduras - refreshed totals
$ws.Cells.Item(103,1).Value = "Honduras"
$ws.Cells.Item(103,2).Value = 139
$ws.Cells.Item(103,3).Value = 29
$ws.Cells.Item(103,4).Value = 3
$ws.Cells.Item(103,5).Value = 129
$ws.Cells.Item(103,6).Value = 4
$ws.Cells.Item(103,7).Value = 4
$ws.Cells.Item(103,8).Value = 7

# Row 136: now "Zambia" (was "Polinesia Francesa") - sheet re-sorted by total cases
$ws.Cells.Item(136,1).Value = "Zambia"
$ws.Cells.Item(136,2).Value = 35
$ws.Cells.Item(136,3).Value = 6
$ws.Cells.Item(136,4).Value = 0
$ws.Cells.Item(136,5).Value = 35
$ws.Cells.Item(136,6).Value = 0
$ws.Cells.Item(136,7).Value = 0
$ws.Cells.Item(136,8).Value = 0

# Row 137: now "Polinesia Francesa" (was "Zambia") - sheet re-sorted by total cases
$ws.Cells.Item(137,1).Value = "Polinesia Francesa"
$ws.Cells.Item(137,2).Value = 35
$ws.Cells.Item(137,3).Value = 5
$ws.Cells.Item(137,4).Value = 0
$ws.Cells.Item(137,5).Value = 35
$ws.Cells.Item(137,6).Value = 2
$ws.Cells.Item(137,7).Value = 0
$ws.Cells.Item(137,8).Value = 0

# Row 143: now "Bermudas" (was "Mali") - sheet re-sorted by total cases
$ws.Cells.Item(143,1).Value = "Bermudas"
$ws.Cells.Item(143,2).Value = 27
$ws.Cells.Item(143,3).Value = 5
$ws.Cells.Item(143,4).Value = 2
$ws.Cells.Item(143,5).Value = 25
$ws.Cells.Item(143,6).Value = 0
$ws.Cells.Item(143,7).Value = 0
$ws.Cells.Item(143,8).Value = 0

# Row 144: now "Mali" (was "Etiopia") - sheet re-sorted by total cases
$ws.Cells.Item(144,1).Value = "Mali"
$ws.Cells.Item(144,2).Value = 25
$ws.Cells.Item(144,3).Value = 7
$ws.Cells.Item(144,4).Value = 0
$ws.Cells.Item(144,5).Value = 23
$ws.Cells.Item(144,6).Value = 0
$ws.Cells.Item(144,7).Value = 1
$ws.Cells.Item(144,8).Value = 2

# Row 145: now "Etiopia" (was "Guinea") - sheet re-sorted by total cases
$ws.Cells.Item(145,1).Value = "Etiopia"
$ws.Cells.Item(145,2).Value = 23
$ws.Cells.Item(145,3).Value = 2
$ws.Cells.Item(145,4).Value = 4
$ws.Cells.Item(145,5).Value = 19
$ws.Cells.Item(145,6).Value = 1
$ws.Cells.Item(145,7).Value = 0
$ws.Cells.Item(145,8).Value = 0

# Row 146: now "Guinea" (was "Bermudas") - sheet re-sorted by total cases
$ws.Cells.Item(146,1).Value = "Guinea"
$ws.Cells.Item(146,2).Value = 22
$ws.Cells.Item(146,3).Value = 6
$ws.Cells.Item(146,4).Value = 0
$ws.Cells.Item(146,5).Value = 22
$ws.Cells.Item(146,6).Value = 0
$ws.Cells.Item(146,7).Value = 0
$ws.Cells.Item(146,8).Value = 0

# Row 169: Siria - refreshed totals
$ws.Cells.Item(169,1).Value = "Siria"
$ws.Cells.Item(169,2).Value = 9
$ws.Cells.Item(169,3).Value = 0
$ws.Cells.Item(169,4).Value = 0
$ws.Cells.Item(169,5).Value = 7
$ws.Cells.Item(169,6).Value = 0
$ws.Cells.Item(169,7).Value = 1
$ws.Cells.Item(169,8).Value = 2

# Row 193: Mauritania - refreshed totals
$ws.Cells.Item(193,1).Value = "Mauritania"
$ws.Cells.Item(193,2).Value = 5
$ws.Cells.Item(193,3).Value = 0
$ws.Cells.Item(193,4).Value = 2
$ws.Cells.Item(193,5).Value = 2
$ws.Cells.Item(193,6).Value = 0
$ws.Cells.Item(193,7).Value = 1
$ws.Cells.Item(193,8).Value = 1
